# Daily attendance processing - 2026-01-20 07:43:33
# Normalize the "Recorded By" column (G) so entries that combine the
# recorder's email with the automated "System" actor are listed with
# "System" first, e.g. "dnasr281@gmail.com, System" -> "System, dnasr281@gmail.com".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

$oldValue = "dnasr281@gmail.com, System"
$newValue = "System, dnasr281@gmail.com"

for ($row = 1; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)   # Column G = "Recorded By"
    if ($cell.Value() -eq $oldValue) {
        $cell.Value = $newValue
    }
}
